# Applies the "cleaned up output" edit to the RMF patch-management report.
#
# The commit condenses/rewrites almost the entire body: dozens of short
# "label: value" / bullet paragraphs are merged into prose paragraphs, bullet
# markers change from "* "/"*** " runs to "- " lists, and whole sections (the
# CVE list, the detailed RMF-controls checklist, the GitLab/Git vulnerability
# breakdown) are replaced by two or three condensed sentences. Because the
# paragraph-level mapping from old to new is not 1:1 (many paragraphs are
# merged, dropped, or reordered), the most faithful way to reproduce the
# diff is to leave the first two paragraphs untouched —
#   1. "Operating System Patch Management RMF Compliance" (Heading2)
#   2. "*** System Overview ***"
# — and wholesale-replace everything after paragraph 2 with the new set of
# paragraphs, in order, reproducing each one's trailing <w:br/> (or lack of
# one, for the very last paragraph) from the target document.

$d = $word.ActiveDocument

# 1. Remove everything after paragraph 2 (keep heading + "System Overview" line).
$anchor = $d.Paragraphs(2)
$tail = $d.Range($anchor.Range.End, $d.Content.End)
$tail.Delete()

# 2. The new body content: each entry is (text, hasTrailingLineBreak).
$items = @(
    ,@("The system is a Debian-based operating system with version #1 SMP PREEMPT_DYNAMIC 6.1.129-1 (2025-03-06). The computer has the name `"kb322-18`" and an IP address of 140.160.138.147.", $true)
    ,@("", $true)
    ,@("*** Patch Status Summary ***", $true)
    ,@("There are several pending updates available:", $true)
    ,@("- Code/stable", $true)
    ,@("- Ure/stable-security", $true)
    ,@("- Git-man/stable-security", $true)
    ,@("- Git/stable-security", $true)
    ,@("", $true)
    ,@("These patches are related to security, specifically addressing vulnerabilities in the git version and Debian's kernel.", $true)
    ,@("", $true)
    ,@("*** Compliance with RMF Controls *** ", $true)
    ,@("In order to comply with the RMF controls, it is recommended that flaw remediation be performed as follows: ", $true)
    ,@("Flaws must be identified through careful monitoring of system logs. Reports on potential issues should be made to a designated individual or team for review and corrective action.", $true)
    ,@("Configuration management should ensure all patches are properly installed and applied consistently across the system. Documentation should also be updated with new patch information.", $true)
    ,@("", $true)
    ,@("Vulnerability checks can be performed through regular scans of the system using tools provided by Debian.", $true)
    ,@("", $true)
    ,@("*** Recommended next steps ***", $true)
    ,@("The recommended next steps are:", $true)
    ,@("- Review and assess updates to ensure compatibility and relevance.", $true)
    ,@("- Schedule patch deployments for all applicable systems in an orderly manner, minimizing disruptions to normal operations.", $true)
    ,@("- Provide guidance for update documentation to facilitate accurate tracking and reporting of patches applied.", $true)
    ,@("", $true)
    ,@("*** Risk Assessment ***", $true)
    ,@("There is a potential risk associated with the pending security updates. The impact level of this vulnerability is high due to its potential to allow unauthorized access or malicious activities on the system. A mitigation plan can be implemented by prioritizing patching of critical systems, monitoring for any suspicious activity, and having an incident response plan in place.", $true)
    ,@("", $true)
    ,@("Please note that without information on the severity of specific CVEs, it's difficult to provide a more detailed risk assessment.", $false)
)

# 3. Re-insert the new paragraphs, one at a time, right after paragraph 2.
#    A trailing line break is produced by appending a manual line break
#    ([char]11, i.e. Shift+Enter) to the text — the engine lowers that to a
#    <w:br/> inside the same run, matching the document's existing
#    "text + <w:br/>" paragraph shape. Paragraph 2 has no explicit style, so
#    paragraphs inserted right after it (and chained from there) naturally
#    stay plain body text too, with no stray pStyle/rsid noise.
$curIndex = 2
foreach ($item in $items) {
    $text = $item[0]
    $hasBreak = $item[1]

    $curPara = $d.Paragraphs($curIndex)
    $curPara.Range.InsertParagraphAfter()
    $curIndex = $curIndex + 1

    $newPara = $d.Paragraphs($curIndex)
    if ($hasBreak) {
        $newPara.Range.Text = $text + [char]11
    } else {
        $newPara.Range.Text = $text
    }
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)

